$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 140, pushing existing rows 140:163 down to 141:164
$ws.Rows.Item(140).Insert()

# Populate the newly inserted row 140 with the new weekly record
$ws.Range("A140").Value = 7
$ws.Range("B140").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C140").Value = "Ñuble"
$ws.Range("D140").Value = 45154
$ws.Range("E140").Value = 16
$ws.Range("F140").Value = 100112031
$ws.Range("G140").Value = "Poroto verde"
$ws.Range("H140").Value = "Magnum"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 50
$ws.Range("K140").Value = 27000
$ws.Range("L140").Value = 27000
$ws.Range("M140").Value = 27000
$ws.Range("N140").Value = "$/malla 25 kilos"
$ws.Range("O140").Value = "Perú"
$ws.Range("P140").Value = 1080
$ws.Range("Q140").Value = 25
$ws.Range("R140").Value = "Hortaliza"
